# Implement parameterization and runmodes for test cases.
#
# Before: a single sheet "test_suite" with columns TCID / Runmode listing the
# test case names (FlightSearchTest, SignInTest) and whether they should run.
#
# After: two new data sheets are introduced that hold the actual
# parameters used by each test case (FlightSearchTest, SignInTest), each
# with its own "runmode" column, plus the original test_suite sheet
# (now listing SignInTest first, then FlightSearchTest).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the three sheets in the right tab order:
#      FlightSearchTest, SignInTest, test_suite
#    The original sheet (currently named "test_suite") is replaced by a
#    fresh sheet of the same name so it starts with no leftover column
#    widths / formatting from the old layout.
# ---------------------------------------------------------------------

$orig = $wb.Worksheets.Item(1)

$signIn = $wb.Worksheets.Add($orig)
$signIn.Name = "SignInTest"

$flightSearch = $wb.Worksheets.Add($wb.Worksheets.Item("SignInTest"))
$flightSearch.Name = "FlightSearchTest"

$origNow = $wb.Worksheets.Item("test_suite")
$suiteNew = $wb.Worksheets.Add($null, $origNow)
$suiteNew.Name = "test_suite_new"

$null = $wb.Worksheets.Item("test_suite").Delete()
$wb.Worksheets.Item("test_suite_new").Name = "test_suite"

# ---------------------------------------------------------------------
# 2. FlightSearchTest sheet: search parameters
# ---------------------------------------------------------------------

$sheet = $wb.Worksheets.Item("FlightSearchTest")

$sheet.Range("A1").Value = "fromCity"
$sheet.Range("B1").Value = "toCity"
$sheet.Range("C1").Value = "departingDate"
$sheet.Range("D1").Value = "returningDate"
$sheet.Range("E1").Value = "adults"
$sheet.Range("F1").Value = "children"
$sheet.Range("G1").Value = "runmode"

$sheet.Range("A2").Value = "Odessa, Ukraine (ODS-Odessa Intl.)"
$sheet.Range("B2").Value = "Paris, France (PAR-All Airports)"
$sheet.Range("C2").Value = "23/06/2018"
$sheet.Range("D2").Value = "31/07/2018"
$sheet.Range("E2").Value = 1
$sheet.Range("F2").Value = 1
$sheet.Range("G2").Value = "Y"

# departingDate / returningDate columns are kept as plain text
$sheet.Columns("C:D").NumberFormat = "@"
$sheet.Columns("C").ColumnWidth = 13.25
$sheet.Columns("D").ColumnWidth = 12.75

$wb.Worksheets.Item("FlightSearchTest").Activate()
$null = $wb.Worksheets.Item("FlightSearchTest").Range("H5").Select()

# ---------------------------------------------------------------------
# 3. SignInTest sheet: credentials (with mailto / hyperlinks on the
#    username column, like the original author's sheet)
# ---------------------------------------------------------------------

$sheet = $wb.Worksheets.Item("SignInTest")

$sheet.Range("A1").Value = "username"
$sheet.Range("B1").Value = "password"
$sheet.Range("C1").Value = "runmode"

$sheet.Range("A2").Value = "test@gmail.com"
$sheet.Range("B2").Value = "passw12rd"
$sheet.Range("C2").Value = "Y"

$sheet.Range("A3").Value = "test2@test"
$sheet.Range("B3").Value = "testpwd"
$sheet.Range("C3").Value = "Y"

$sheet.Hyperlinks.Add($sheet.Range("A2"), "mailto:test@gmail.com")
$sheet.Hyperlinks.Add($sheet.Range("A3"), "mailto:test2@test")

$sheet.Columns("A").ColumnWidth = 13.42
$sheet.Columns("B").ColumnWidth = 13.59

$wb.Worksheets.Item("SignInTest").Activate()
$null = $wb.Worksheets.Item("SignInTest").Range("E5").Select()

# ---------------------------------------------------------------------
# 4. test_suite sheet: which test cases run (SignInTest now listed
#    before FlightSearchTest)
# ---------------------------------------------------------------------

$sheet = $wb.Worksheets.Item("test_suite")

$sheet.Range("A1").Value = "TCID"
$sheet.Range("B1").Value = "Runmode"

$sheet.Range("A2").Value = "SignInTest"
$sheet.Range("B2").Value = "Y"

$sheet.Range("A3").Value = "FlightSearchTest"
$sheet.Range("B3").Value = "Y"

$wb.Worksheets.Item("test_suite").Activate()
$null = $wb.Worksheets.Item("test_suite").Range("A2").Select()

# ---------------------------------------------------------------------
# 5. Final UI state: SignInTest is the active/selected tab.
# ---------------------------------------------------------------------

$wb.Worksheets.Item("SignInTest").Activate()

Write-Host "done"
